# daily auto push: 2026-01-16 18:48 UTC
# Two new log rows (2026/01/16 22:xx and 2026/01/17 02:xx) were recorded and
# inserted into the daily log sheet right before the 2026/12/29 entries,
# pushing every following row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 664/665 (existing rows 664.. shift down to 666..)
$ws.Range("A664:A665").EntireRow.Insert()

# Row 664: 2026/01/16 金 22:xx -> rank 201
$ws.Range("A664").Value = "'2026/01/16"
$ws.Range("B664").Value = "'金"
$ws.Range("C664").Value = 22
$ws.Range("D664").Value = 201

# Row 665: 2026/01/17 土 02:xx -> rank 201
$ws.Range("A665").Value = "'2026/01/17"
$ws.Range("B665").Value = "'土"
$ws.Range("C665").Value = 2
$ws.Range("D665").Value = 201
